$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 2.08
$ws.Range("J4").Value = 4.4
$ws.Range("AC4").Value = 11
$ws.Range("T5").Value = 2.1
$ws.Range("X5").Value = 21
$ws.Range("F6").Value = 5.2
$ws.Range("Q6").Value = 2.08
$ws.Range("T6").Value = 1.99
$ws.Range("F7").Value = 1.22
$ws.Range("G7").Value = 1.24
$ws.Range("J7").Value = 7.8
$ws.Range("O7").Value = 1.14
$ws.Range("T7").Value = 1.99
$ws.Range("U7").Value = 1.89
$ws.Range("X7").Value = 100
$ws.Range("Y7").Value = 320
$ws.Range("AB7").Value = 15.5
$ws.Range("AC7").Value = 19.5
$ws.Range("AD7").Value = 230
$ws.Range("AE7").Value = 230
$ws.Range("AH7").Value = 75
$ws.Range("AM7").Value = 160
$ws.Range("AO7").Value = 250
$ws.Range("N8").Value = 5.4
$ws.Range("U8").Value = 2.44
$ws.Range("F9").Value = 1.58
$ws.Range("J9").Value = 4.5
$ws.Range("K9").Value = 4.8
$ws.Range("Q9").Value = 1.69
$ws.Range("K10").Value = 5.1
$ws.Range("P11").Value = 1.87
$ws.Range("F12").Value = 1.6
$ws.Range("G12").Value = 1.62
$ws.Range("H12").Value = 5.8
$ws.Range("I12").Value = 6.4
$ws.Range("J12").Value = 4.6
$ws.Range("T12").Value = 1.8
$ws.Range("U12").Value = 2.18
$ws.Range("X12").Value = 26
$ws.Range("Z12").Value = 140
$ws.Range("AB12").Value = 10
$ws.Range("AD12").Value = 29
$ws.Range("AI12").Value = 990
$ws.Range("AJ12").Value = 15.5
$ws.Range("AN12").Value = 7.4
$ws.Range("G14").Value = 2.52
$ws.Range("J14").Value = 3.15
$ws.Range("P14").Value = 1.74
$ws.Range("F16").Value = 1.63
$ws.Range("J16").Value = 3.85
$ws.Range("F17").Value = 2.46
$ws.Range("G17").Value = 2.76
$ws.Range("H17").Value = 2.78
$ws.Range("I17").Value = 3.1
$ws.Range("J17").Value = 3.4
$ws.Range("K17").Value = 3.7
$ws.Range("F19").Value = 2.56
$ws.Range("J19").Value = 3.3
$ws.Range("K19").Value = 3.55
$ws.Range("P19").Value = 1.82
$ws.Range("Q19").Value = 2.14
$ws.Range("Q20").Value = 1.92
$ws.Range("H22").Value = 1.52
$ws.Range("P22").Value = 2.48
$ws.Range("AA23").Value = 410
$ws.Range("AD23").Value = 1000
$ws.Range("AH23").Value = 1000
$ws.Range("AI23").Value = 150
$ws.Range("F25").Value = 1.78
$ws.Range("Q25").Value = 1.94
$ws.Range("AA26").Value = 19.5
$ws.Range("K27").Value = 6.2
$ws.Range("K28").Value = 5
$ws.Range("AD28").Value = 44
$ws.Range("AE28").Value = 130
$ws.Range("AL28").Value = 38
$ws.Range("G29").Value = 1.24
$ws.Range("H29").Value = 17
$ws.Range("J29").Value = 7
$ws.Range("P29").Value = 2.4
$ws.Range("S29").Value = 2.66
$ws.Range("U29").Value = 1.65
$ws.Range("X29").Value = 27
$ws.Range("AB29").Value = 8.800000000000001
$ws.Range("P31").Value = 1.99
$ws.Range("Q31").Value = 1.8
$ws.Range("H32").Value = 1.64
$ws.Range("J33").Value = 4.2
$ws.Range("P33").Value = 2.16
$ws.Range("Q33").Value = 1.69
$ws.Range("F34").Value = 4.5
$ws.Range("H34").Value = 1.88
$ws.Range("K34").Value = 3.95
$ws.Range("Q34").Value = 1.83
$ws.Range("F36").Value = 6.6
$ws.Range("H36").Value = 1.45
$ws.Range("F37").Value = 1.69
$ws.Range("J37").Value = 3.8
$ws.Range("H38").Value = 2.28
$ws.Range("I38").Value = 2.46
$ws.Range("F39").Value = 2.66
$ws.Range("H39").Value = 2.72
$ws.Range("I39").Value = 3.2
$ws.Range("J39").Value = 3.1
$ws.Range("P39").Value = 1.7
$ws.Range("Q39").Value = 2.34
$ws.Range("G40").Value = 1.97
$ws.Range("H40").Value = 5.2
